$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (prices + 1h volume deltas in columns D/E).
# Some new price strings parse as plain numbers under Excel's normal type
# inference, which would silently drop significant trailing zeros (e.g.
# "0.610" -> 0.61); force those specific cells to Text first so the
# literal digit string is preserved exactly, matching the source data.

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "41.358.66"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "2.179.40"
$ws.Range("E3").Value = "  -1.58%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "237.87"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("D6").Value = "0.610"
$ws.Range("D7").Value = "70.22"
$ws.Range("E7").Value = "  -3.96%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.576"
$ws.Range("E9").Value = "  -4.88%  "
$ws.Range("D10").Value = "39.64"
$ws.Range("E10").Value = "  -7.52%  "
$ws.Range("D11").Value = "0.0921"
$ws.Range("E11").Value = "  -3.26%  "
$ws.Range("D12").Value = "54.69"
$ws.Range("E12").Value = "  -4.80%  "
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("D14").Value = "6.76"
$ws.Range("E14").Value = "  -4.53%  "
$ws.Range("D15").Value = "2.498.97"
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("D16").Value = "14.37"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").Value = "0.796"
$ws.Range("E17").Value = "  -4.70%  "
$ws.Range("D18").Value = "2.159.35"
$ws.Range("E18").Value = "  -2.20%  "
$ws.Range("D19").Value = "41.151.74"
$ws.Range("E19").Value = "  -1.51%  "
$ws.Range("E20").Value = "  -6.95%  "
$ws.Range("D21").Value = "70.64"
$ws.Range("E21").Value = "  -3.43%  "
$ws.Range("D22").Value = "5.88"
$ws.Range("E22").Value = "  -4.35%  "
$ws.Range("D23").Value = "226.57"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").Value = "9.41"
$ws.Range("E24").Value = "  -8.55%  "
$ws.Range("D25").Value = "1.91"
$ws.Range("E25").Value = "  -8.44%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "10.79"
$ws.Range("E27").Value = "  -7.32%  "
$ws.Range("E28").Value = "  -3.35%  "
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  -2.11%  "
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("D31").Value = "167.68"
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("D32").Value = "19.94"
$ws.Range("E32").Value = "  -3.09%  "
$ws.Range("D33").Value = "30.35"
$ws.Range("E33").Value = "  +4.23%  "
$ws.Range("D34").Value = "0.0764"
$ws.Range("E34").Value = "  -3.85%  "
$ws.Range("D35").Value = "5.13"
$ws.Range("E35").Value = "  -9.80%  "
$ws.Range("D36").Value = "0.121"
$ws.Range("E36").Value = "  -3.17%  "
$ws.Range("D37").Value = "0.102"
$ws.Range("E37").Value = "  -7.74%  "
$ws.Range("D38").Value = "4.10"
$ws.Range("E38").Value = "  -3.39%  "
$ws.Range("D39").Value = "0.0283"
$ws.Range("E39").Value = "  -5.52%  "
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D41").Value = "11.86"
$ws.Range("E41").Value = "  -11.98%  "
$ws.Range("D42").Value = "5.39"
$ws.Range("E42").Value = "  -3.85%  "
$ws.Range("D43").Value = "59.06"
$ws.Range("E43").Value = "  -10.83%  "
$ws.Range("D44").Value = "0.191"
$ws.Range("E44").Value = "  -3.43%  "
$ws.Range("D45").Value = "8.31"
$ws.Range("E45").Value = "  -3.96%  "
$ws.Range("D46").Value = "0.0968"
$ws.Range("E46").Value = "  -3.44%  "
$ws.Range("D47").Value = "97.60"
$ws.Range("E47").Value = "  -5.75%  "
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("D49").Value = "1.13"
$ws.Range("E49").Value = "  -3.16%  "
$ws.Range("D50").Value = "2.20"
$ws.Range("E50").Value = "  -7.97%  "
$ws.Range("D51").Value = "2.61"
$ws.Range("E51").Value = "  -2.67%  "
